$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Merge "ICES division" (col H) and "ICES subdivision" (col I) into a single
# "ICES area" column. Column I is removed; everything to its right shifts left.
$ws.Columns("I").Delete()
$ws.Range("H1").Value = "ICES area"
